# Flaming Fox review: move/rework the "Meta description" blurb.
#
# 1. The "Meta description" paragraph (bold label + blurb) right after the
#    H1 title is removed entirely.
# 2. A new bold paragraph carrying the page title text
#    ("Play Flaming Fox Free: Review of Unique Board Slot") is inserted
#    just before the final (italic "Prompt: ...") paragraph.
# 3. That final paragraph's italic text is swapped from the old image-prompt
#    text to the blurb text that used to live in the meta description
#    paragraph.

$d = $word.ActiveDocument

# --- Step 1: drop the "Meta description" paragraph (2nd paragraph, right
# after the H1 title) -- Range.Delete() removes the paragraph mark too, so
# the following paragraph collapses up cleanly. ----------------------------
$metaLabel = "Meta description"
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($metaLabel)) {
        $targetParaIndex = $i
        break
    }
}
if ($targetParaIndex -ge 1) {
    $d.Paragraphs.Item($targetParaIndex).Range.Delete() | Out-Null
}

# --- Step 2: insert a new bold paragraph with the title text right before
# the last paragraph in the document. --------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Flaming Fox Free: Review of Unique Board Slot</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($titleXml) | Out-Null

# --- Step 3: swap the closing "Prompt: ..." paragraph's text for the old
# meta-description blurb (keeps its existing italic run formatting). ------
$oldPrompt = "Prompt: Create a cartoon-style feature image for Flaming Fox that features a happy Maya warrior with glasses. The image should prominently display the Flaming Fox title and include elements of the Chinese temple and ninja fox theme of the game, such as flames and oriental architecture. The Maya warrior should be depicted wearing glasses and holding a winning combination of ninja swords surrounded by flames. The background should be vibrant and colourful, incorporating elements of Chinese culture and martial arts, as well as the signature flaming fox theme of the game. The image should be eye-catching and appealing to players who enjoy action-packed slot games with a lot of personality and unique features."
$newBlurb = "Experience the unique gameplay and bonus features of Flaming Fox for free. Review of this high-volatility board slot with high rewards and intricate graphics."

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newBlurb, 2) | Out-Null
